# "novos ajustes a pedido do usuario"
#
# The workbook has 3 sheets:
#   1. B.I. (ajuste -)
#   2. B.I. (ajuste +)
#   3. B.I. (intercompany)
#
# Target state has only 2 sheets:
#   1. B.I. Intercompany   (was "B.I. (intercompany)", moved to front)
#   2. B.I. Passivo        (was "B.I. (ajuste -)")
#
# i.e. "B.I. (ajuste +)" is removed entirely, "B.I. (intercompany)" is
# renamed and promoted to the first tab, and "B.I. (ajuste -)" is renamed
# to "B.I. Passivo" and kept as the (now) second, active tab.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Drop the "ajuste +" sheet completely.
$wb.Worksheets.Item("B.I. (ajuste +)").Delete() | Out-Null

# Move the intercompany sheet to the very first tab position.
$wb.Worksheets.Item("B.I. (intercompany)").Move($wb.Worksheets.Item(1))

# Rename the remaining sheets to their new names.
$wb.Worksheets.Item("B.I. (intercompany)").Name = "B.I. Intercompany"
$wb.Worksheets.Item("B.I. (ajuste -)").Name = "B.I. Passivo"

# "B.I. Passivo" remains the active/selected tab.
$wb.Worksheets.Item("B.I. Passivo").Activate()

$excel.DisplayAlerts = $true
